$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the blank/duplicate serial number in C4 (was "58A022747")
$ws.Range("C4").Value = ""

# Fill in the blank "x" in the PRINT DD-1750 column for the row that now
# carries the previously-duplicated serial number (D5)
$ws.Range("D5").Value = "x"

# Refresh the active cell selection, matching Excel's own post-edit state
$ws.Range("E11").Select()
